# Scheduled-runner market data refresh for Hades_Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for a
# handful of leve rows across the crafting-profession sheets, reflecting
# fresh market-board prices. A few rows gain/lose an HQ profit cell (M)
# because the item now does/doesn't have an HQ price on the board.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70: Consecrating Congregation
$ws.Range("H70").Value = 791.2727
$ws.Range("I70").Value = 595.4
$ws.Range("K70").Value = 1786.2
$ws.Range("M70").Value = -1516.2

# Row 73: Curbing the Contagion (L)
$ws.Range("H73").Value = 791.2727
$ws.Range("I73").Value = 595.4
$ws.Range("K73").Value = 1786.2
$ws.Range("M73").Value = -850.1999999999998

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1668948.6
$ws.Range("I137").Value = 2633094
$ws.Range("J137").Value = 3606.5
$ws.Range("K137").Value = 7899282
$ws.Range("L137").Value = 10819.5
$ws.Range("M137").Value = -7896732
$ws.Range("N137").Value = -15919.5

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3706803.5
$ws.Range("I138").Value = 5614.1665
$ws.Range("J138").Value = 4169452.2
$ws.Range("K138").Value = 16842.4995
$ws.Range("L138").Value = 12508356.6
$ws.Range("M138").Value = -11702.4995
$ws.Range("N138").Value = -12518636.6

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 23461.625
$ws.Range("I32").Value = 22566.875
$ws.Range("J32").Value = 26145.875
$ws.Range("K32").Value = 22566.875
$ws.Range("L32").Value = 26145.875
$ws.Range("M32").Value = -22279.875
$ws.Range("N32").Value = -26719.875

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 1001.7143
$ws.Range("I45").Value = 918.6667
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 918.6667
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -541.6667
$ws.Range("N45").Value = -2254

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1309.9584
$ws.Range("I110").Value = 890.5333000000001
$ws.Range("K110").Value = 890.5333000000001
$ws.Range("M110").Value = 1154.4667

$ws = $wb.Worksheets.Item("BSM")
# Row 25: Tools of the Trade
$ws.Range("H25").Value = 1575.5

# Row 94: High Steal
$ws.Range("H94").Value = 428.4762
$ws.Range("I94").Value = 384.6316
$ws.Range("K94").Value = 384.6316
$ws.Range("M94").Value = 66.36840000000001

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2264.6482
$ws.Range("I134").Value = 2250.7737
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 6752.321100000001
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -4217.321100000001
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 234567.4
$ws.Range("I31").Value = 39982.117
$ws.Range("J31").Value = 429152.7
$ws.Range("K31").Value = 39982.117
$ws.Range("L31").Value = 429152.7
$ws.Range("M31").Value = -39687.117
$ws.Range("N31").Value = -429742.7

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 234567.4
$ws.Range("I34").Value = 39982.117
$ws.Range("J34").Value = 429152.7
$ws.Range("K34").Value = 39982.117
$ws.Range("L34").Value = 429152.7
$ws.Range("M34").Value = -39780.117
$ws.Range("N34").Value = -429556.7

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 3091.647
$ws.Range("I62").Value = 2782
$ws.Range("J62").Value = 3220.6667
$ws.Range("K62").Value = 2782
$ws.Range("L62").Value = 3220.6667
$ws.Range("M62").Value = -2158
$ws.Range("N62").Value = -4468.6667

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 3091.647
$ws.Range("I65").Value = 2782
$ws.Range("J65").Value = 3220.6667
$ws.Range("K65").Value = 13910
$ws.Range("L65").Value = 16103.3335
$ws.Range("M65").Value = -10790
$ws.Range("N65").Value = -22343.3335

# Row 86: Birch, Please
$ws.Range("H86").Value = 3213.5
$ws.Range("I86").Value = 2900
$ws.Range("J86").Value = 3258.2856
$ws.Range("K86").Value = 2900
$ws.Range("L86").Value = 3258.2856
$ws.Range("M86").Value = -1777
$ws.Range("N86").Value = -5504.2856

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 3213.5
$ws.Range("I89").Value = 2900
$ws.Range("J89").Value = 3258.2856
$ws.Range("K89").Value = 14500
$ws.Range("L89").Value = 16291.428
$ws.Range("M89").Value = -8884
$ws.Range("N89").Value = -27523.428

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water
$ws.Range("H4").Value = 6669192
$ws.Range("J4").Value = 7695033
$ws.Range("L4").Value = 23085099
$ws.Range("N4").Value = -23085323

# Row 108: Meet for Meat
$ws.Range("H108").Value = 405
$ws.Range("I108").Value = 405
$ws.Range("K108").Value = 1215
$ws.Range("M108").Value = 1665

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 616.4583
$ws.Range("I113").Value = 561.6875
$ws.Range("J113").Value = 726
$ws.Range("K113").Value = 1685.0625
$ws.Range("L113").Value = 2178
$ws.Range("M113").Value = 484.9375
$ws.Range("N113").Value = -6518

# Row 114: One Last Meal
$ws.Range("H114").Value = 10101527
$ws.Range("J114").Value = 15873566
$ws.Range("L114").Value = 47620698
$ws.Range("N114").Value = -47627206

# Row 117: A Good Omen
$ws.Range("H117").Value = 4167653.2
$ws.Range("I117").Value = 636.8570999999999
$ws.Range("J117").Value = 7408666
$ws.Range("K117").Value = 1910.5713
$ws.Range("L117").Value = 22225998
$ws.Range("M117").Value = 1531.4287
$ws.Range("N117").Value = -22232882

# Row 118: Teetotally
$ws.Range("H118").Value = 3265.5625
$ws.Range("I118").Value = 641.5
$ws.Range("J118").Value = 4840
$ws.Range("K118").Value = 1924.5
$ws.Range("L118").Value = 14520
$ws.Range("M118").Value = -681.5
$ws.Range("N118").Value = -17006

# Row 129: Comfort Food
$ws.Range("H129").Value = 4388055
$ws.Range("I129").Value = 1256.5555
$ws.Range("J129").Value = 8336173.5
$ws.Range("K129").Value = 3769.6665
$ws.Range("L129").Value = 25008520.5
$ws.Range("M129").Value = 1230.3335
$ws.Range("N129").Value = -25018520.5

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 1847
$ws.Range("I137").Value = 873.3333
$ws.Range("J137").Value = 3599.6
$ws.Range("K137").Value = 2619.9999
$ws.Range("L137").Value = 10798.8
$ws.Range("M137").Value = 2480.0001
$ws.Range("N137").Value = -20998.8

# Row 141: Ocean Explosion
$ws.Range("H141").Value = 9741.1875
$ws.Range("I141").Value = 3705.9
$ws.Range("K141").Value = 11117.7
$ws.Range("M141").Value = -5937.700000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 2187.7
$ws.Range("I82").Value = 1377.6
$ws.Range("J82").Value = 2997.8
$ws.Range("K82").Value = 1377.6
$ws.Range("L82").Value = 2997.8
$ws.Range("M82").Value = -1016.6
$ws.Range("N82").Value = -3719.8

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 2187.7
$ws.Range("I85").Value = 1377.6
$ws.Range("J85").Value = 2997.8
$ws.Range("K85").Value = 1377.6
$ws.Range("L85").Value = 2997.8
$ws.Range("M85").Value = -129.5999999999999
$ws.Range("N85").Value = -5493.8

# Row 122: Hell on Leather
$ws.Range("H122").Value = 3372.5
$ws.Range("I122").Value = 3518.3333
$ws.Range("K122").Value = 10554.9999
$ws.Range("M122").Value = -8104.999899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 1377
$ws.Range("I122").Value = 1377
$ws.Range("K122").Value = 4131
$ws.Range("M122").Value = -1681

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 55369.71
$ws.Range("I136").Value = 41585.56
$ws.Range("K136").Value = 124756.68
$ws.Range("M136").Value = -122206.68
